$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2 = @(128, 0, 128, 29.784663655714184, 58.927911801228497, -36.487081171480732)
    3 = @(0, 0, 255, 32.297009439844494, 79.187517397197226, -107.86016288933186)
    4 = @(0, 128, 0, 46.227430171917774, -51.698494452473625, 49.896839611104149)
    5 = @(255, 255, 0, 97.139263431696349, -21.553728492530066, 94.477963319459832)
    6 = @(255, 192, 151, 82.484871067825324, 17.53855336966803, 29.45464984614874)
    7 = @(255, 0, 0, 53.240788867616104, 80.092494286414734, 67.20319139735453)
    8 = @(255, 0, 0, 53.240788867616104, 80.092494286414734, 67.20319139735453)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $vals[$col - 1]
    }
}
